$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =====================================================================
# 1) Overwrite the data values for the existing rows 3-7 (clusters 1-5)
#    with the re-drafted ("redactada nuevamente") numbers, and add the
#    three new rows 8-10 (clusters 6-8) for the FRE clustering proposal.
# =====================================================================

# --- Row 3 (cluster 1) ---
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 475.5
$ws.Cells.Item(3, 3).Value = 71.447478785000001
$ws.Cells.Item(3, 4).Value = 91
$ws.Cells.Item(3, 5).Value = 0.76470588235294101
$ws.Cells.Item(3, 6).Value = 0.88743369638835001
$ws.Cells.Item(3, 7).Value = 0.75
$ws.Cells.Item(3, 8).Value = 0.41666666666666602
$ws.Cells.Item(3, 9).Value = 7.6541100000000002

# --- Row 4 (cluster 2) ---
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 2714
$ws.Cells.Item(4, 3).Value = 978.38835269000003
$ws.Cells.Item(4, 4).Value = 1054
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.96164976306739802
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 1
$ws.Cells.Item(4, 9).Value = 0.081818181818181804

# --- Row 5 (cluster 3) ---
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 821.93333333333305
$ws.Cells.Item(5, 3).Value = 80.272733311636003
$ws.Cells.Item(5, 4).Value = 149.933333333333
$ws.Cells.Item(5, 5).Value = 0.91764705882352904
$ws.Cells.Item(5, 6).Value = 0.96146745545502599
$ws.Cells.Item(5, 7).Value = 1
$ws.Cells.Item(5, 8).Value = 0.54814740740740697
$ws.Cells.Item(5, 9).Value = 0.91869512254752705

# --- Row 6 (cluster 4) ---
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 1252
$ws.Cells.Item(6, 3).Value = 15.5245885625
$ws.Cells.Item(6, 4).Value = 90.5
$ws.Cells.Item(6, 5).Value = 0.94117647058823495
$ws.Cells.Item(6, 6).Value = 0.80235460422413496
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = 0.47222222222222199
$ws.Cells.Item(6, 9).Value = 100

# --- Row 7 (cluster 5) ---
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 615.33333333333303
$ws.Cells.Item(7, 3).Value = 57.380750808333296
$ws.Cells.Item(7, 4).Value = 101.333333333333
$ws.Cells.Item(7, 5).Value = 0.88235294117647001
$ws.Cells.Item(7, 6).Value = 0.209039759485211
$ws.Cells.Item(7, 7).Value = 0.95833333333333304
$ws.Cells.Item(7, 8).Value = 0.74074074074074003
$ws.Cells.Item(7, 9).Value = 0.48735538461538402

# --- Row 8 (cluster 6, new) ---
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 102
$ws.Cells.Item(8, 3).Value = 38.902520000000003
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 0.58823529411764697
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 0.875
$ws.Cells.Item(8, 8).Value = 0.88888888888888795
$ws.Cells.Item(8, 9).Value = 90

# --- Row 9 (cluster 7, new) ---
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 201.4
$ws.Cells.Item(9, 3).Value = 19.011901713127997
$ws.Cells.Item(9, 4).Value = 12.4
$ws.Cells.Item(9, 5).Value = 0.494117647058823
$ws.Cells.Item(9, 6).Value = 0.91935327750763396
$ws.Cells.Item(9, 7).Value = 1
$ws.Cells.Item(9, 8).Value = 0.4
$ws.Cells.Item(9, 9).Value = 3.4593068545454502

# --- Row 10 (cluster 8, new) ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 1667
$ws.Cells.Item(10, 3).Value = 276.62430657599998
$ws.Cells.Item(10, 4).Value = 453
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.87146318071491902
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0.17674000000000001

# =====================================================================
# 2) Formatting. Apply alignment first, then number formats, so the
#    engine reuses/merges style records instead of fragmenting them.
# =====================================================================

# --- Rows 3-7 keep the original centered + vertically centered look ---
$ws.Range("A3:I7").HorizontalAlignment = -4108
$ws.Range("A3:I7").VerticalAlignment = -4108

$ws.Range("B3:D7").NumberFormat = "0"
$ws.Range("E3:E7").NumberFormat = "0.00"
$ws.Range("G3:H7").NumberFormat = "0.00"
$ws.Range("F3:F7").NumberFormat = "0.000"
$ws.Range("I3:I7").NumberFormat = "0.000"

# --- New rows 8-10 use horizontal centering only (no vertical centering) ---
$ws.Range("A8:I10").HorizontalAlignment = -4108

$ws.Range("B8:D10").NumberFormat = "0"
$ws.Range("E8:E10").NumberFormat = "0.00"
$ws.Range("G8:H10").NumberFormat = "0.00"
$ws.Range("F8:F10").NumberFormat = "0.000"
$ws.Range("I8:I10").NumberFormat = "0.000"

# =====================================================================
# 3) Extend the per-column color-scale conditional formatting so it
#    covers the new rows (was 3:7, now 3:10) for every metric column.
# =====================================================================
$cfCols = @("B", "C", "D", "E", "F", "G", "H", "I")
foreach ($col in $cfCols) {
    $oldRange = $ws.Range($col + "3:" + $col + "7")
    $newRange = $ws.Range($col + "3:" + $col + "10")
    $fcCount = $oldRange.FormatConditions.Count
    for ($i = 1; $i -le $fcCount; $i++) {
        $oldRange.FormatConditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# =====================================================================
# 4) Move the active selection to G3 (matches the saved cursor state).
# =====================================================================
$ws.Range("G3").Select()
